# Auto-generated script to update cryptos worksheet data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = "'" + $value
    $ws.Cells.Item($row, $col).Style = "Normal"
}

Set-TextCell 2 4 "29.549.49"
Set-TextCell 2 5 "  -0.80%  "
Set-TextCell 3 4 "1.852.15"
Set-TextCell 3 5 "  -0.37%  "
Set-TextCell 4 5 "  -0.13%  "
Set-TextCell 5 4 "243.45"
Set-TextCell 6 5 "  -1.00%  "
Set-TextCell 8 2 "Cardano"
Set-TextCell 8 3 "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextCell 8 4 "0.3004"
Set-TextCell 8 5 "  -0.20%  "
Set-TextCell 9 2 "Dogecoin"
Set-TextCell 9 3 "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextCell 9 4 "0.07478"
Set-TextCell 9 5 "  -0.68%  "
Set-TextCell 10 2 "Solana"
Set-TextCell 10 3 "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextCell 10 4 "24.29"
Set-TextCell 10 5 "  +0.34%  "
Set-TextCell 11 2 "TRON"
Set-TextCell 11 3 "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell 11 4 "0.07627"
Set-TextCell 11 5 "  -0.90%  "
Set-TextCell 12 2 "WrappedEther"
Set-TextCell 12 3 "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell 12 4 "1.902.34"
Set-TextCell 12 5 "  +2.16%  "
Set-TextCell 13 2 "Polkadot"
Set-TextCell 13 3 "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell 13 4 "5.032"
Set-TextCell 13 5 "  -0.69%  "
Set-TextCell 14 2 "Polygon"
Set-TextCell 14 3 "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextCell 14 4 "0.6861"
Set-TextCell 14 5 "  +0.08%  "
Set-TextCell 15 2 "Litecoin"
Set-TextCell 15 3 "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell 15 4 "83.60"
Set-TextCell 15 5 "  -0.61%  "
Set-TextCell 16 2 "ShibaInu"
Set-TextCell 16 3 "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell 16 4 "0.000009545"
Set-TextCell 16 5 "  +0.71%  "
Set-TextCell 17 2 "Uniswap"
Set-TextCell 17 3 "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell 17 4 "6.169"
Set-TextCell 17 5 "  +1.22%  "
Set-TextCell 18 2 "WrappedBTC"
Set-TextCell 18 3 "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell 18 4 "29.557.29"
Set-TextCell 18 5 "  -0.76%  "
Set-TextCell 19 2 "WrappedliquidstakedEther2.0"
Set-TextCell 19 3 "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextCell 19 4 "2.119.36"
Set-TextCell 19 5 "  -0.22%  "
Set-TextCell 20 2 "BitcoinCash"
Set-TextCell 20 3 "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell 20 4 "236.29"
Set-TextCell 20 5 "  -2.11%  "
Set-TextCell 21 2 "Avalanche"
Set-TextCell 21 3 "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell 21 4 "12.57"
Set-TextCell 21 5 "  -1.12%  "
Set-TextCell 22 2 "Dai"
Set-TextCell 22 3 "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell 22 4 "1.000"
Set-TextCell 22 5 "  -0.04%  "
Set-TextCell 23 2 "Chainlink"
Set-TextCell 23 3 "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell 23 4 "7.726"
Set-TextCell 23 5 "  +3.86%  "
Set-TextCell 24 2 "BinanceUSD"
Set-TextCell 24 3 "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextCell 24 4 "1.000"
Set-TextCell 24 5 "  -0.17%  "
Set-TextCell 25 2 "Monero"
Set-TextCell 25 3 "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell 25 4 "157.32"
Set-TextCell 25 5 "  -1.14%  "
Set-TextCell 26 2 "Stellar"
Set-TextCell 26 3 "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell 26 4 "0.1403"
Set-TextCell 26 5 "  -1.95%  "
Set-TextCell 27 2 "Cosmos"
Set-TextCell 27 3 "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell 27 4 "8.508"
Set-TextCell 27 5 "  -0.57%  "
Set-TextCell 28 2 "EthereumClassic"
Set-TextCell 28 3 "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell 28 4 "17.77"
Set-TextCell 28 5 "  -1.26%  "
Set-TextCell 29 2 "PancakeSwap"
Set-TextCell 29 3 "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell 29 4 "1.490"
Set-TextCell 29 5 "  -0.93%  "
Set-TextCell 30 2 "Hedera"
Set-TextCell 30 3 "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell 30 4 "0.05994"
Set-TextCell 30 5 "  -2.19%  "
Set-TextCell 31 2 "Toncoin"
Set-TextCell 31 3 "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell 31 4 "1.250"
Set-TextCell 31 5 "  -1.91%  "
Set-TextCell 32 2 "Filecoin"
Set-TextCell 32 3 "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell 32 4 "4.122"
Set-TextCell 32 5 "  -1.04%  "
Set-TextCell 33 2 "InternetComputer(DFINITY)"
Set-TextCell 33 3 "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell 33 4 "4.074"
Set-TextCell 33 5 "  -1.35%  "
Set-TextCell 34 2 "LidoDAOToken"
Set-TextCell 34 3 "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell 34 4 "1.872"
Set-TextCell 34 5 "  -0.76%  "
Set-TextCell 35 2 "ARBITRUM"
Set-TextCell 35 3 "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell 35 4 "1.177"
Set-TextCell 35 5 "  +1.57%  "
Set-TextCell 36 2 "ImmutableX"
Set-TextCell 36 3 "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell 36 4 "0.7207"
Set-TextCell 36 5 "  -1.94%  "
Set-TextCell 37 2 "HuobiToken"
Set-TextCell 37 3 "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell 37 4 "2.603"
Set-TextCell 37 5 "  -0.25%  "
Set-TextCell 38 2 "MXToken"
Set-TextCell 38 3 "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell 38 4 "2.802"
Set-TextCell 38 5 "  -1.93%  "
Set-TextCell 39 2 "VeChain"
Set-TextCell 39 3 "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell 39 4 "0.01776"
Set-TextCell 39 5 "  -1.37%  "
Set-TextCell 40 2 "Maker"
Set-TextCell 40 3 "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell 40 4 "1.203.51"
Set-TextCell 40 5 "  -1.45%  "
Set-TextCell 41 2 "TrustWalletToken"
Set-TextCell 41 3 "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell 41 4 "0.9096"
Set-TextCell 41 5 "  -2.49%  "
Set-TextCell 42 2 "FraxShare"
Set-TextCell 42 3 "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell 42 4 "6.164"
Set-TextCell 42 5 "  -1.73%  "
Set-TextCell 43 2 "RocketPoolETH"
Set-TextCell 43 3 "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextCell 43 4 "2.043.07"
Set-TextCell 43 5 "  +0.46%  "
Set-TextCell 44 2 "PaxDollar"
Set-TextCell 44 3 "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextCell 44 4 "0.9994"
Set-TextCell 44 5 "  -0.21%  "
Set-TextCell 45 2 "Quant"
Set-TextCell 45 3 "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell 45 4 "101.95"
Set-TextCell 45 5 "  -0.14%  "
Set-TextCell 46 2 "Aave"
Set-TextCell 46 3 "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell 46 4 "66.76"
Set-TextCell 46 5 "  +0.26%  "
Set-TextCell 47 2 "Aptos"
Set-TextCell 47 3 "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell 47 4 "7.295"
Set-TextCell 47 5 "  +8.65%  "
Set-TextCell 48 2 "BabyDogeCoin"
Set-TextCell 48 3 "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell 48 4 "0.00000000118"
Set-TextCell 48 5 "  -4.01%  "
Set-TextCell 49 2 "TheSandbox"
Set-TextCell 49 3 "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextCell 49 4 "0.4035"
Set-TextCell 49 5 "  -1.39%  "
Set-TextCell 50 2 "EnergySwap"
Set-TextCell 50 3 "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell 50 4 "9.113"
Set-TextCell 50 5 "  -2.56%  "
Set-TextCell 51 2 "RenderToken"
Set-TextCell 51 3 "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell 51 4 "1.659"
Set-TextCell 51 5 "  +0.93%  "

Write-Host "Applied $(183) cell updates"
